$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Overview": swap the two data rows (bfc.. now first, eb53.. now
# second) and mark both as "Handed back: in sync with en-US" (the
# "Ready for handoff" status disappears entirely).
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = "bfc56466-f424-4c71-a9e5-4b645e843490.md"
$wsOverview.Range("B2").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("D2").Value2 = "2016-03-24 00:51:35"

$wsOverview.Range("A3").Value2 = "eb5361da-e598-49e1-a781-2298c27002b2.md"
$wsOverview.Range("B3").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value2 = "Handed back: in sync with en-US"
$wsOverview.Range("D3").Value2 = "2016-03-24 00:50:08"

# Hyperlinks keep pointing at the very same targets they always did
# (rId2 keeps resolving to eb5361da's github url, rId3 keeps resolving
# to bfc56466's github url) -- only the *displayed* text is swapped to
# track the new row contents.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a26c2ce4ffc3aeae213fd111388237465e00a0c4/e2e/eb5361da-e598-49e1-a781-2298c27002b2.md",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.md"
) | Out-Null
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f2aca258f9e1e0f31e6c090b2896561b71794e03/e2e/bfc56466-f424-4c71-a9e5-4b645e843490.md",
    "",
    "",
    "eb5361da-e598-49e1-a781-2298c27002b2.md"
) | Out-Null

# -----------------------------------------------------------------
# Sheet "zh-cn": same row-content swap as Overview, plus the handback
# datetime (col H) for the bfc56466 file advances to 00:51:56.
# -----------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value2 = "bfc56466-f424-4c71-a9e5-4b645e843490.md"
$wsZh.Range("B2").Value2 = ".md"
$wsZh.Range("C2").Value2 = "Handed back: in sync with en-US"
$wsZh.Range("D2").Value2 = "bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf"
$wsZh.Range("E2").Value2 = "2016-03-24 00:51:31"
$wsZh.Range("F2").Value2 = "bfc56466-f424-4c71-a9e5-4b645e843490.md"
$wsZh.Range("G2").Value2 = "bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf"
$wsZh.Range("H2").Value2 = "2016-03-24 00:51:56"
$wsZh.Range("J2").Value2 = "Include"

$wsZh.Range("A3").Value2 = "eb5361da-e598-49e1-a781-2298c27002b2.md"
$wsZh.Range("B3").Value2 = ".md"
$wsZh.Range("C3").Value2 = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value2 = "eb5361da-e598-49e1-a781-2298c27002b2.7b12573b727babc04b54bfddf3d680bab711f64c.zh-cn.xlf"
$wsZh.Range("E3").Value2 = "2016-03-24 00:50:02"
$wsZh.Range("F3").Value2 = "eb5361da-e598-49e1-a781-2298c27002b2.md"
$wsZh.Range("G3").Value2 = "eb5361da-e598-49e1-a781-2298c27002b2.7b12573b727babc04b54bfddf3d680bab711f64c.zh-cn.xlf"
$wsZh.Range("H3").Value2 = "2016-03-24 00:50:38"
$wsZh.Range("J3").Value2 = "Include"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a26c2ce4ffc3aeae213fd111388237465e00a0c4/e2e/eb5361da-e598-49e1-a781-2298c27002b2.md",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e4c1973a36e434833cbfef11cc03e03f64d4723/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/eb5361da-e598-49e1-a781-2298c27002b2.7b12573b727babc04b54bfddf3d680bab711f64c.zh-cn.xlf",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/d7fdb2afcc9853fd64c2276c9f4a8d3dcc892ad1/e2e/eb5361da-e598-49e1-a781-2298c27002b2.md",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b56bc8424907eb88b1710c83cfb009349da617a9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/eb5361da-e598-49e1-a781-2298c27002b2.7b12573b727babc04b54bfddf3d680bab711f64c.zh-cn.xlf",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f2aca258f9e1e0f31e6c090b2896561b71794e03/e2e/bfc56466-f424-4c71-a9e5-4b645e843490.md",
    "",
    "",
    "eb5361da-e598-49e1-a781-2298c27002b2.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d334afd6da1c26c86212a43bf72de19f6bb40da/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf",
    "",
    "",
    "eb5361da-e598-49e1-a781-2298c27002b2.7b12573b727babc04b54bfddf3d680bab711f64c.zh-cn.xlf"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/d7fdb2afcc9853fd64c2276c9f4a8d3dcc892ad1/e2e/bfc56466-f424-4c71-a9e5-4b645e843490.md",
    "",
    "",
    "eb5361da-e598-49e1-a781-2298c27002b2.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b56bc8424907eb88b1710c83cfb009349da617a9/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.zh-cn.xlf",
    "",
    "",
    "eb5361da-e598-49e1-a781-2298c27002b2.7b12573b727babc04b54bfddf3d680bab711f64c.zh-cn.xlf"
) | Out-Null

# -----------------------------------------------------------------
# Sheet "de-de": same row-content swap, handback datetime (col H) for
# bfc56466 advances to 00:52:03.
# -----------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value2 = "bfc56466-f424-4c71-a9e5-4b645e843490.md"
$wsDe.Range("B2").Value2 = ".md"
$wsDe.Range("C2").Value2 = "Handed back: in sync with en-US"
$wsDe.Range("D2").Value2 = "bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf"
$wsDe.Range("E2").Value2 = "2016-03-24 00:51:35"
$wsDe.Range("F2").Value2 = "bfc56466-f424-4c71-a9e5-4b645e843490.md"
$wsDe.Range("G2").Value2 = "bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf"
$wsDe.Range("H2").Value2 = "2016-03-24 00:52:03"
$wsDe.Range("J2").Value2 = "Include"

$wsDe.Range("A3").Value2 = "eb5361da-e598-49e1-a781-2298c27002b2.md"
$wsDe.Range("B3").Value2 = ".md"
$wsDe.Range("C3").Value2 = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value2 = "eb5361da-e598-49e1-a781-2298c27002b2.7b12573b727babc04b54bfddf3d680bab711f64c.de-de.xlf"
$wsDe.Range("E3").Value2 = "2016-03-24 00:50:08"
$wsDe.Range("F3").Value2 = "eb5361da-e598-49e1-a781-2298c27002b2.md"
$wsDe.Range("G3").Value2 = "eb5361da-e598-49e1-a781-2298c27002b2.7b12573b727babc04b54bfddf3d680bab711f64c.de-de.xlf"
$wsDe.Range("H3").Value2 = "2016-03-24 00:50:47"
$wsDe.Range("J3").Value2 = "Include"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a26c2ce4ffc3aeae213fd111388237465e00a0c4/e2e/eb5361da-e598-49e1-a781-2298c27002b2.md",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6d843090e2ed1158392d1a82d7f435837cdb140e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/eb5361da-e598-49e1-a781-2298c27002b2.7b12573b727babc04b54bfddf3d680bab711f64c.de-de.xlf",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/dd1589d6978e304c1bc5b6e8d22689ff9506c99b/e2e/eb5361da-e598-49e1-a781-2298c27002b2.md",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a0146951b1da12cc2312ce7341409b59eb213d58/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/eb5361da-e598-49e1-a781-2298c27002b2.7b12573b727babc04b54bfddf3d680bab711f64c.de-de.xlf",
    "",
    "",
    "bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f2aca258f9e1e0f31e6c090b2896561b71794e03/e2e/bfc56466-f424-4c71-a9e5-4b645e843490.md",
    "",
    "",
    "eb5361da-e598-49e1-a781-2298c27002b2.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b6f8c878dd73f33eff0fc18d2062985d2eca290/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf",
    "",
    "",
    "eb5361da-e598-49e1-a781-2298c27002b2.7b12573b727babc04b54bfddf3d680bab711f64c.de-de.xlf"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/dd1589d6978e304c1bc5b6e8d22689ff9506c99b/e2e/bfc56466-f424-4c71-a9e5-4b645e843490.md",
    "",
    "",
    "eb5361da-e598-49e1-a781-2298c27002b2.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a0146951b1da12cc2312ce7341409b59eb213d58/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bfc56466-f424-4c71-a9e5-4b645e843490.fbb322ddca8ef794e052c77312888adbd6407e5c.de-de.xlf",
    "",
    "",
    "eb5361da-e598-49e1-a781-2298c27002b2.7b12573b727babc04b54bfddf3d680bab711f64c.de-de.xlf"
) | Out-Null

Write-Output "Generate Report for Handback: done"
